$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Roll the quarterly window forward one quarter -------------------------
# Drop the oldest quarter's column (D = "فصل سوم منتهی به 1399/06") and
# append a new column at the end (M) for the newest quarter
# ("فصل اول منتهی به 1401/12"), shifting every other quarter's data one
# column to the left. This mirrors the periodic "roll the data window"
# update described in the commit message.
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).ColumnWidth = 28.17

# --- Row 8: quarter-end headers --------------------------------------------
# Columns D..L already hold the correct (shifted) labels; only the new
# column M needs the brand new quarter label.
$ws.Range("M8").Value = "فصل اول منتهی به 1401/12"

# --- Row 9: publication dates ----------------------------------------------
# After the shift, I9/L9 still hold stale labels that need correcting, and
# M9 needs the new publication date.
$ws.Range("I9").Value = "1402-01-30 (2)"
$ws.Range("L9").Value = "1402-01-30 (3)"
$ws.Range("M9").Value = "1402-01-30"

# --- Financial data: new quarter's figures (column M) ----------------------
$ws.Range("M11").Value = 3139808
$ws.Range("M12").Value = -1637053
$ws.Range("M13").Value = 1502755
$ws.Range("M14").Value = -135259
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("M17").Value = 1367496
$ws.Range("M18").Value = 0
$ws.Range("M19").Value = 307792
$ws.Range("M20").Value = 1675288
$ws.Range("M21").Value = -177940
$ws.Range("M22").Value = 1497348
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 1497348
$ws.Range("M25").Value = 607
$ws.Range("M26").Value = 2466000
$ws.Range("M27").Value = 607
